$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H39").Value = 2281.6667
$ws.Range("I39").Value = 42
$ws.Range("J39").Value = 4073.4
$ws.Range("K39").Value = 126
$ws.Range("L39").Value = 12220.2
$ws.Range("M39").Value = 170
$ws.Range("N39").Value = -12812.2
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H74").Value = 20006040
$ws.Range("I74").Value = 29414588
$ws.Range("J74").Value = 12874.375
$ws.Range("K74").Value = 29414588
$ws.Range("L74").Value = 12874.375
$ws.Range("M74").Value = -29413652
$ws.Range("N74").Value = -14746.375
$ws.Range("H77").Value = 20006040
$ws.Range("I77").Value = 29414588
$ws.Range("J77").Value = 12874.375
$ws.Range("K77").Value = 147072940
$ws.Range("L77").Value = 64371.875
$ws.Range("M77").Value = -147068260
$ws.Range("N77").Value = -73731.875
$ws.Range("H96").Value = 643.3570999999999
$ws.Range("I96").Value = 539.38464
$ws.Range("J96").Value = 1995
$ws.Range("K96").Value = 1618.15392
$ws.Range("L96").Value = 5985
$ws.Range("M96").Value = -245.15392
$ws.Range("N96").Value = -8731
$ws.Range("H98").Value = 4447.2
$ws.Range("I98").Value = 4447.2
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 4447.2
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -2949.2
$ws.Range("H111").Value = 10418743
$ws.Range("I111").Value = 25000680
$ws.Range("J111").Value = 3074.2856
$ws.Range("K111").Value = 75002040
$ws.Range("L111").Value = 9222.856800000001
$ws.Range("M111").Value = -74998973
$ws.Range("N111").Value = -15356.8568
$ws.Range("H122").Value = 4447.2
$ws.Range("I122").Value = 4447.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13341.6
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -10891.6
$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 50000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H138").Value = 1566291.4
$ws.Range("I138").Value = 1803.6538
$ws.Range("J138").Value = 2636730.2
$ws.Range("K138").Value = 5410.9614
$ws.Range("L138").Value = 7910190.600000001
$ws.Range("M138").Value = -270.9614000000001
$ws.Range("N138").Value = -7920470.600000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 433.66666
$ws.Range("I5").Value = 300.5
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 300.5
$ws.Range("L5").Value = 700
$ws.Range("M5").Value = -188.5
$ws.Range("N5").Value = -924
$ws.Range("H32").Value = 4768880.5
$ws.Range("I32").Value = 5133012.5
$ws.Range("J32").Value = 35166.332
$ws.Range("K32").Value = 5133012.5
$ws.Range("L32").Value = 35166.332
$ws.Range("M32").Value = -5132725.5
$ws.Range("N32").Value = -35740.332
$ws.Range("H45").Value = 2540.95
$ws.Range("I45").Value = 1939.5555
$ws.Range("J45").Value = 3033
$ws.Range("K45").Value = 1939.5555
$ws.Range("L45").Value = 3033
$ws.Range("M45").Value = -1562.5555
$ws.Range("N45").Value = -3787
$ws.Range("H63").Value = 1975
$ws.Range("I63").Value = 1983.3334
$ws.Range("J63").Value = 1958.3334
$ws.Range("K63").Value = 1983.3334
$ws.Range("L63").Value = 1958.3334
$ws.Range("M63").Value = -1297.3334
$ws.Range("N63").Value = -3330.3334
$ws.Range("H66").Value = 1975
$ws.Range("I66").Value = 1983.3334
$ws.Range("J66").Value = 1958.3334
$ws.Range("K66").Value = 9916.666999999999
$ws.Range("L66").Value = 9791.666999999999
$ws.Range("M66").Value = -6484.666999999999
$ws.Range("N66").Value = -16655.667
$ws.Range("H74").Value = 32925.09
$ws.Range("I74").Value = 43238.125
$ws.Range("J74").Value = 5423.6665
$ws.Range("K74").Value = 43238.125
$ws.Range("L74").Value = 5423.6665
$ws.Range("M74").Value = -42364.125
$ws.Range("N74").Value = -7171.6665
$ws.Range("H77").Value = 32925.09
$ws.Range("I77").Value = 43238.125
$ws.Range("J77").Value = 5423.6665
$ws.Range("K77").Value = 216190.625
$ws.Range("L77").Value = 27118.3325
$ws.Range("M77").Value = -211822.625
$ws.Range("N77").Value = -35854.3325
$ws.Range("H97").Value = 3970727.8
$ws.Range("I97").Value = 2074.1428
$ws.Range("J97").Value = 11908035
$ws.Range("K97").Value = 2074.1428
$ws.Range("L97").Value = 11908035
$ws.Range("M97").Value = -1578.1428
$ws.Range("N97").Value = -11909027
$ws.Range("H122").Value = 3826.5
$ws.Range("I122").Value = 1547.5834
$ws.Range("J122").Value = 5779.857
$ws.Range("K122").Value = 4642.7502
$ws.Range("L122").Value = 17339.571
$ws.Range("M122").Value = -2192.7502
$ws.Range("N122").Value = -22239.571
$ws.Range("H132").Value = 5423.985
$ws.Range("I132").Value = 3922.3958
$ws.Range("J132").Value = 9428.223
$ws.Range("K132").Value = 11767.1874
$ws.Range("L132").Value = 28284.669
$ws.Range("M132").Value = -9237.187399999999
$ws.Range("N132").Value = -33344.669
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 433.66666
$ws.Range("I4").Value = 300.5
$ws.Range("J4").Value = 700
$ws.Range("K4").Value = 300.5
$ws.Range("L4").Value = 700
$ws.Range("M4").Value = -185.5
$ws.Range("N4").Value = -930
$ws.Range("H11").Value = 1873.6666
$ws.Range("I11").Value = 310.5
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 310.5
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = -170.5
$ws.Range("N11").Value = -5280
$ws.Range("H107").Value = 38797560
$ws.Range("I107").Value = 41671240
$ws.Range("J107").Value = 2845
$ws.Range("K107").Value = 41671240
$ws.Range("L107").Value = 2845
$ws.Range("M107").Value = -41669320
$ws.Range("N107").Value = -6685
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 2449.0667
$ws.Range("I29").Value = 2019
$ws.Range("J29").Value = 2479.7856
$ws.Range("K29").Value = 2019
$ws.Range("L29").Value = 2479.7856
$ws.Range("M29").Value = -1726
$ws.Range("N29").Value = -3065.7856
$ws.Range("H31").Value = 6035.8394
$ws.Range("I31").Value = 1625.6538
$ws.Range("J31").Value = 9858
$ws.Range("K31").Value = 1625.6538
$ws.Range("L31").Value = 9858
$ws.Range("M31").Value = -1330.6538
$ws.Range("N31").Value = -10448
$ws.Range("H34").Value = 6035.8394
$ws.Range("I34").Value = 1625.6538
$ws.Range("J34").Value = 9858
$ws.Range("K34").Value = 1625.6538
$ws.Range("L34").Value = 9858
$ws.Range("M34").Value = -1423.6538
$ws.Range("N34").Value = -10262
$ws.Range("H122").Value = 2024
$ws.Range("I122").Value = 1426.8334
$ws.Range("J122").Value = 3047.7144
$ws.Range("K122").Value = 4280.5002
$ws.Range("L122").Value = 9143.143199999999
$ws.Range("M122").Value = -1830.5002
$ws.Range("N122").Value = -14043.1432
$ws.Range("H132").Value = 5754.927
$ws.Range("I132").Value = 3601.913
$ws.Range("J132").Value = 8506
$ws.Range("K132").Value = 10805.739
$ws.Range("L132").Value = 25518
$ws.Range("M132").Value = -8275.739
$ws.Range("N132").Value = -30578
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1167.5
$ws.Range("I16").Value = 723.3333
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 2169.9999
$ws.Range("L16").Value = 7500
$ws.Range("M16").Value = -1996.9999
$ws.Range("N16").Value = -7846
$ws.Range("H19").Value = 110
$ws.Range("I19").Value = 110
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 330
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -156
$ws.Range("N19").ClearContents()
$ws.Range("H32").Value = 142857230
$ws.Range("I32").Value = 500000060
$ws.Range("J32").Value = 98
$ws.Range("K32").Value = 1500000180
$ws.Range("L32").Value = 294
$ws.Range("M32").Value = -1499999897
$ws.Range("N32").Value = -860
$ws.Range("H50").Value = 111112450
$ws.Range("I50").Value = 111112450
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 333337350
$ws.Range("L50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -333336869
$ws.Range("H53").Value = 111112450
$ws.Range("I53").Value = 111112450
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 333337350
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -333336869
$ws.Range("H55").Value = 4173895.8
$ws.Range("I55").Value = 504
$ws.Range("J55").Value = 4355347.5
$ws.Range("K55").Value = 1512
$ws.Range("L55").Value = 13066042.5
$ws.Range("M55").Value = -1335
$ws.Range("N55").Value = -13066396.5
$ws.Range("H58").Value = 2676.25
$ws.Range("I58").Value = 2901.6667
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 8705.000100000001
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -8577.000100000001
$ws.Range("N58").Value = -6256
$ws.Range("H121").Value = 25000748
$ws.Range("I121").Value = 50000250
$ws.Range("J121").Value = 1248.5
$ws.Range("K121").Value = 150000750
$ws.Range("L121").Value = 3745.5
$ws.Range("M121").Value = -149999440
$ws.Range("N121").Value = -6365.5
$ws.Range("H132").Value = 5218.6387
$ws.Range("I132").Value = 2469.077
$ws.Range("J132").Value = 6772.7393
$ws.Range("K132").Value = 22221.693
$ws.Range("L132").Value = 60954.6537
$ws.Range("M132").Value = -19691.693
$ws.Range("N132").Value = -66014.6537
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 73636.86
$ws.Range("I80").Value = 1981.8572
$ws.Range("J80").Value = 145291.86
$ws.Range("K80").Value = 1981.8572
$ws.Range("L80").Value = 145291.86
$ws.Range("M80").Value = -983.8571999999999
$ws.Range("N80").Value = -147287.86
$ws.Range("H83").Value = 73636.86
$ws.Range("I83").Value = 1981.8572
$ws.Range("J83").Value = 145291.86
$ws.Range("K83").Value = 9909.286
$ws.Range("L83").Value = 726459.2999999999
$ws.Range("M83").Value = -4917.286
$ws.Range("N83").Value = -736443.2999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6747.222
$ws.Range("I132").Value = 3507.353
$ws.Range("J132").Value = 8714.286
$ws.Range("K132").Value = 10522.059
$ws.Range("L132").Value = 26142.858
$ws.Range("M132").Value = -7992.059000000001
$ws.Range("N132").Value = -31202.858
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 9804718
$ws.Range("I107").Value = 382.26315
$ws.Range("J107").Value = 22223542
$ws.Range("K107").Value = 1146.78945
$ws.Range("L107").Value = 66670626
$ws.Range("M107").Value = 773.21055
$ws.Range("N107").Value = -66674466
$ws.Range("H122").Value = 2591.0732
$ws.Range("I122").Value = 2018.2413
$ws.Range("J122").Value = 3975.4167
$ws.Range("K122").Value = 6054.7239
$ws.Range("L122").Value = 11926.2501
$ws.Range("M122").Value = -3604.7239
$ws.Range("N122").Value = -16826.2501
$ws.Range("H132").Value = 2818.5
$ws.Range("I132").Value = 1227.5
$ws.Range("J132").Value = 6000.5
$ws.Range("K132").Value = 3682.5
$ws.Range("L132").Value = 18001.5
$ws.Range("M132").Value = -1152.5
$ws.Range("N132").Value = -23061.5
$ws.Range("H139").Value = 90181.5
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 90181.5
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 90181.5
$ws.Range("N139").Value = -100461.5
